$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 188.21053
$ws.Range("I11").Value = 188.21053
$ws.Range("K11").Value = 188.21053
$ws.Range("M11").Value = -48.21053000000001
$ws.Range("H17").Value = 387631.22
$ws.Range("J17").Value = 387631.22
$ws.Range("L17").Value = 1162893.66
$ws.Range("N17").Value = -1163229.66
$ws.Range("H38").Value = 514.2353
$ws.Range("I38").Value = 66.3
$ws.Range("J38").Value = 1154.1428
$ws.Range("K38").Value = 198.9
$ws.Range("L38").Value = 3462.4284
$ws.Range("M38").Value = 173.1
$ws.Range("N38").Value = -4206.428400000001
$ws.Range("H42").Value = 95.14286
$ws.Range("I42").Value = 25.7
$ws.Range("J42").Value = 268.75
$ws.Range("K42").Value = 77.1
$ws.Range("L42").Value = 806.25
$ws.Range("M42").Value = 152.9
$ws.Range("N42").Value = -1266.25
$ws.Range("H58").Value = 330.23077
$ws.Range("I58").Value = 299.45456
$ws.Range("J58").Value = 499.5
$ws.Range("K58").Value = 898.36368
$ws.Range("L58").Value = 1498.5
$ws.Range("M58").Value = -748.36368
$ws.Range("N58").Value = -1798.5
$ws.Range("H86").Value = 7803.3125
$ws.Range("I86").Value = 6874.778
$ws.Range("J86").Value = 8997.143
$ws.Range("K86").Value = 6874.778
$ws.Range("L86").Value = 8997.143
$ws.Range("M86").Value = -5751.778
$ws.Range("N86").Value = -11243.143
$ws.Range("H89").Value = 7803.3125
$ws.Range("I89").Value = 6874.778
$ws.Range("J89").Value = 8997.143
$ws.Range("K89").Value = 34373.89
$ws.Range("L89").Value = 44985.715
$ws.Range("M89").Value = -28757.89
$ws.Range("N89").Value = -56217.715
$ws.Range("H98").Value = 1155.6
$ws.Range("I98").Value = 1047.3914
$ws.Range("J98").Value = 2400.0
$ws.Range("K98").Value = 1047.3914
$ws.Range("L98").Value = 2400.0
$ws.Range("M98").Value = 450.6086
$ws.Range("N98").Value = -5396.0
$ws.Range("H100").Value = 4305.5
$ws.Range("I100").Value = 4305.5
$ws.Range("J100").Value = 0.0
$ws.Range("K100").Value = 4305.5
$ws.Range("L100").Value = 0.0
$ws.Range("M100").Value = -3764.5
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 1155.6
$ws.Range("I122").Value = 1047.3914
$ws.Range("J122").Value = 2400.0
$ws.Range("K122").Value = 3142.1742
$ws.Range("L122").Value = 7200.0
$ws.Range("M122").Value = -692.1741999999999
$ws.Range("N122").Value = -12100.0
$ws.Range("H132").Value = 5848893.5
$ws.Range("I132").Value = 6803679.0
$ws.Range("K132").Value = 20411037.0
$ws.Range("M132").Value = -20408507.0
$ws.Range("H137").Value = 1712.9667
$ws.Range("I137").Value = 1583.5
$ws.Range("J137").Value = 2230.8333
$ws.Range("K137").Value = 4750.5
$ws.Range("L137").Value = 6692.499899999999
$ws.Range("M137").Value = -2200.5
$ws.Range("N137").Value = -11792.4999
$ws.Range("H141").Value = 2220.45
$ws.Range("I141").Value = 1788.7646
$ws.Range("J141").Value = 4666.6665
$ws.Range("K141").Value = 5366.293799999999
$ws.Range("L141").Value = 13999.9995
$ws.Range("M141").Value = -186.2937999999995
$ws.Range("N141").Value = -24359.9995
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 30429.0
$ws.Range("I26").Value = 8248.75
$ws.Range("K26").Value = 8248.75
$ws.Range("M26").Value = -7918.75
$ws.Range("H32").Value = 2676.3135
$ws.Range("I32").Value = 1895.6508
$ws.Range("J32").Value = 14971.75
$ws.Range("K32").Value = 1895.6508
$ws.Range("L32").Value = 14971.75
$ws.Range("M32").Value = -1608.6508
$ws.Range("N32").Value = -15545.75
$ws.Range("H45").Value = 5916.3335
$ws.Range("I45").Value = 5981.7856
$ws.Range("K45").Value = 5981.7856
$ws.Range("M45").Value = -5604.7856
$ws.Range("H74").Value = 11264.263
$ws.Range("I74").Value = 1827.6666
$ws.Range("J74").Value = 27441.285
$ws.Range("K74").Value = 1827.6666
$ws.Range("L74").Value = 27441.285
$ws.Range("M74").Value = -953.6666
$ws.Range("N74").Value = -29189.285
$ws.Range("H77").Value = 11264.263
$ws.Range("I77").Value = 1827.6666
$ws.Range("J77").Value = 27441.285
$ws.Range("K77").Value = 9138.333
$ws.Range("L77").Value = 137206.425
$ws.Range("M77").Value = -4770.333000000001
$ws.Range("N77").Value = -145942.425
$ws.Range("H122").Value = 1262.4546
$ws.Range("I122").Value = 1262.4546
$ws.Range("K122").Value = 3787.3638
$ws.Range("M122").Value = -1337.3638
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 812.0
$ws.Range("I64").Value = 1110.25
$ws.Range("J64").Value = 613.1667
$ws.Range("K64").Value = 1110.25
$ws.Range("L64").Value = 613.1667
$ws.Range("M64").Value = -885.25
$ws.Range("N64").Value = -1063.1667
$ws.Range("H67").Value = 812.0
$ws.Range("I67").Value = 1110.25
$ws.Range("J67").Value = 613.1667
$ws.Range("K67").Value = 1110.25
$ws.Range("L67").Value = 613.1667
$ws.Range("M67").Value = -330.25
$ws.Range("N67").Value = -2173.1667
$ws.Range("H86").Value = 4199.8
$ws.Range("J86").Value = 7501.75
$ws.Range("L86").Value = 7501.75
$ws.Range("N86").Value = -9747.75
$ws.Range("H89").Value = 4199.8
$ws.Range("J89").Value = 7501.75
$ws.Range("L89").Value = 37508.75
$ws.Range("N89").Value = -48740.75
$ws.Range("H105").Value = 2224.6155
$ws.Range("I105").Value = 2112.0
$ws.Range("J105").Value = 2600.0
$ws.Range("K105").Value = 2112.0
$ws.Range("L105").Value = 2600.0
$ws.Range("M105").Value = -365.0
$ws.Range("N105").Value = -6094.0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1604.9231
$ws.Range("I12").Value = 705.6667
$ws.Range("J12").Value = 2375.7144
$ws.Range("K12").Value = 705.6667
$ws.Range("L12").Value = 2375.7144
$ws.Range("M12").Value = -535.6667
$ws.Range("N12").Value = -2715.7144
$ws.Range("H31").Value = 41988.652
$ws.Range("I31").Value = 48906.477
$ws.Range("J31").Value = 12933.8
$ws.Range("K31").Value = 48906.477
$ws.Range("L31").Value = 12933.8
$ws.Range("M31").Value = -48611.477
$ws.Range("N31").Value = -13523.8
$ws.Range("H34").Value = 41988.652
$ws.Range("I34").Value = 48906.477
$ws.Range("J34").Value = 12933.8
$ws.Range("K34").Value = 48906.477
$ws.Range("L34").Value = 12933.8
$ws.Range("M34").Value = -48704.477
$ws.Range("N34").Value = -13337.8
$ws.Range("H41").Value = 105.0
$ws.Range("H62").Value = 4077.5
$ws.Range("I62").Value = 2796.8
$ws.Range("J62").Value = 5358.2
$ws.Range("K62").Value = 2796.8
$ws.Range("L62").Value = 5358.2
$ws.Range("M62").Value = -2172.8
$ws.Range("N62").Value = -6606.2
$ws.Range("H65").Value = 4077.5
$ws.Range("I65").Value = 2796.8
$ws.Range("J65").Value = 5358.2
$ws.Range("K65").Value = 13984.0
$ws.Range("L65").Value = 26791.0
$ws.Range("M65").Value = -10864.0
$ws.Range("N65").Value = -33031.0
$ws.Range("H103").Value = 10000.0
$ws.Range("I103").Value = 10000.0
$ws.Range("K103").Value = 10000.0
$ws.Range("M103").Value = -8828.0
$ws.Range("H122").Value = 729.0
$ws.Range("I122").Value = 752.04346
$ws.Range("K122").Value = 2256.13038
$ws.Range("M122").Value = 193.8696199999999
$ws.Range("H132").Value = 3329.525
$ws.Range("I132").Value = 3337.6943
$ws.Range("J132").Value = 3256.0
$ws.Range("K132").Value = 10013.0829
$ws.Range("L132").Value = 9768.0
$ws.Range("M132").Value = -7483.082900000001
$ws.Range("N132").Value = -14828.0
$ws.Range("H134").Value = 6556.17
$ws.Range("I134").Value = 4405.59
$ws.Range("K134").Value = 13216.77
$ws.Range("M134").Value = -10681.77
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1200.0
$ws.Range("J132").Value = 1666.6666
$ws.Range("L132").Value = 14999.9994
$ws.Range("N132").Value = -20059.9994
$ws.Range("H134").Value = 4267.5
$ws.Range("I134").Value = 2591.4285
$ws.Range("K134").Value = 7774.2855
$ws.Range("M134").Value = -2704.2855
$ws.Range("H137").Value = 2796.353
$ws.Range("I137").Value = 1331.4546
$ws.Range("K137").Value = 3994.3638
$ws.Range("M137").Value = 1105.6362
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 2461.8
$ws.Range("I22").Value = 800.0
$ws.Range("J22").Value = 2877.25
$ws.Range("K22").Value = 800.0
$ws.Range("L22").Value = 2877.25
$ws.Range("N22").Value = -3935.25
$ws.Range("M22").Value = -271.0
$ws.Range("H102").Value = 974.4375
$ws.Range("I102").Value = 974.4375
$ws.Range("K102").Value = 974.4375
$ws.Range("M102").Value = 647.5625
$ws.Range("H126").Value = 13673.467
$ws.Range("I126").Value = 17163.818
$ws.Range("K126").Value = 51491.454
$ws.Range("M126").Value = -49021.454
$ws.Range("H132").Value = 2368.0
$ws.Range("I132").Value = 2196.4443
$ws.Range("K132").Value = 6589.3329
$ws.Range("M132").Value = -4059.3329
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 19233356.0
$ws.Range("J16").Value = 4771.3335
$ws.Range("L16").Value = 4771.3335
$ws.Range("N16").Value = -5111.3335
$ws.Range("H46").Value = 2405.6428
$ws.Range("I46").Value = 1279.8
$ws.Range("K46").Value = 1279.8
$ws.Range("M46").Value = -1091.8
$ws.Range("H55").Value = 572.2
$ws.Range("I55").Value = 567.63635
$ws.Range("K55").Value = 567.63635
$ws.Range("M55").Value = -394.63635
$ws.Range("H93").Value = 6799.8
$ws.Range("I93").Value = 6799.8
$ws.Range("K93").Value = 6799.8
$ws.Range("M93").Value = -5551.8
$ws.Range("H132").Value = 3674.182
$ws.Range("I132").Value = 3250.6667
$ws.Range("K132").Value = 9752.000100000001
$ws.Range("M132").Value = -7222.000100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0.0
$ws.Range("J48").Value = 0.0
$ws.Range("L48").Value = 0.0
$ws.Range("N48").ClearContents()
$ws.Range("H112").Value = 42833.168
$ws.Range("J112").Value = 42833.168
$ws.Range("L112").Value = 42833.168
$ws.Range("N112").Value = -45787.168
$ws.Range("H126").Value = 2414.3914
$ws.Range("I126").Value = 2283.625
$ws.Range("J126").Value = 2713.2856
$ws.Range("K126").Value = 6850.875
$ws.Range("L126").Value = 8139.8568
$ws.Range("M126").Value = -4380.875
$ws.Range("N126").Value = -13079.8568
